$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$excel.ScreenUpdating = $false

# 1. Fix D274 value (20206613 -> 20206619)
$ws.Cells.Item(274, 4).Value = 20206619

# 2. Seed formatting for the new rows 328:352 by copying row 327's
#    formats. Columns A:L always carry data, so copy those in one shot.
#    Column M/N are mutually exclusive per-row (only one of the two is
#    ever populated, matching the source survey branching), so format only
#    the single M-or-N cell that will actually hold a value per row -
#    leaving the other column cell completely absent, same as the source.
$ws.Range("A327:L327").Copy()
$ws.Range("A328:L352").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Rows("328:352").RowHeight = 15.75

$ws.Range("N327:N327").Copy()
$ws.Range("N328:N328").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("N327:N327").Copy()
$ws.Range("M329:M329").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("N327:N327").Copy()
$ws.Range("N330:N330").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("N327:N327").Copy()
$ws.Range("N331:N331").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("N327:N327").Copy()
$ws.Range("N332:N332").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("N327:N327").Copy()
$ws.Range("N333:N333").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("N327:N327").Copy()
$ws.Range("M334:M334").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("N327:N327").Copy()
$ws.Range("M335:M335").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("N327:N327").Copy()
$ws.Range("M336:M336").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("N327:N327").Copy()
$ws.Range("N337:N337").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("N327:N327").Copy()
$ws.Range("N338:N338").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("N327:N327").Copy()
$ws.Range("M339:M339").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("N327:N327").Copy()
$ws.Range("N340:N340").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("N327:N327").Copy()
$ws.Range("M341:M341").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("N327:N327").Copy()
$ws.Range("N342:N342").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("N327:N327").Copy()
$ws.Range("M343:M343").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("N327:N327").Copy()
$ws.Range("N344:N344").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("N327:N327").Copy()
$ws.Range("N345:N345").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("N327:N327").Copy()
$ws.Range("N346:N346").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("N327:N327").Copy()
$ws.Range("N347:N347").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("N327:N327").Copy()
$ws.Range("M348:M348").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("N327:N327").Copy()
$ws.Range("M349:M349").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("N327:N327").Copy()
$ws.Range("N350:N350").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("N327:N327").Copy()
$ws.Range("M351:M351").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("N327:N327").Copy()
$ws.Range("M352:M352").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# 3. Fill in new row data (rows 328-352)
# Row 328
$ws.Cells.Item(328, 1).Value = 45192.674386851853
$ws.Cells.Item(328, 2).Value = "gys5785@naver.com"
$ws.Cells.Item(328, 3).Value = "정치행정학과"
$ws.Cells.Item(328, 4).Value = 20232402
$ws.Cells.Item(328, 5).Value = "고형승"
$ws.Cells.Item(328, 6).Value = "‘조(租)’는 공전(公田)의 경작자가 국고에 상납하는 지대 또는 사전(私田)의 경작자가 전주에게 바치는 지대를 뜻한다."
$ws.Cells.Item(328, 7).Value = 0.9
$ws.Cells.Item(328, 8).Value = "4:6"
$ws.Cells.Item(328, 9).Value = "15분의 1"
$ws.Cells.Item(328, 10).Value = "130만호, 5백만명"
$ws.Cells.Item(328, 11).Value = "평안"
$ws.Cells.Item(328, 12).Value = "Black"
$ws.Cells.Item(328, 14).Value = "노동자가 과도한 연장근로를 받을 수 있어 반대한다."

# Row 329
$ws.Cells.Item(329, 1).Value = 45192.680134490744
$ws.Cells.Item(329, 2).Value = "minjoo902@naver.com"
$ws.Cells.Item(329, 3).Value = "금융재무학과"
$ws.Cells.Item(329, 4).Value = 20192827
$ws.Cells.Item(329, 5).Value = "김민주"
$ws.Cells.Item(329, 6).Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Cells.Item(329, 7).Value = 0.1
$ws.Cells.Item(329, 8).Value = "6:4"
$ws.Cells.Item(329, 9).Value = "20분의 1"
$ws.Cells.Item(329, 10).Value = "20만호, 69만명"
$ws.Cells.Item(329, 11).Value = "충청"
$ws.Cells.Item(329, 12).Value = "Red"
$ws.Cells.Item(329, 13).Value = "모름/무응답"

# Row 330
$ws.Cells.Item(330, 1).Value = 45192.6832569213
$ws.Cells.Item(330, 2).Value = "tkdgjs9768@naver.com"
$ws.Cells.Item(330, 3).Value = "경제학과"
$ws.Cells.Item(330, 4).Value = 20212837
$ws.Cells.Item(330, 5).Value = "임상헌"
$ws.Cells.Item(330, 6).Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Cells.Item(330, 7).Value = 0.1
$ws.Cells.Item(330, 8).Value = "6:4"
$ws.Cells.Item(330, 9).Value = "20분의 1"
$ws.Cells.Item(330, 10).Value = "20만호, 69만명"
$ws.Cells.Item(330, 11).Value = "충청"
$ws.Cells.Item(330, 12).Value = "Black"
$ws.Cells.Item(330, 14).Value = "찬성한다."

# Row 331
$ws.Cells.Item(331, 1).Value = 45192.684675416662
$ws.Cells.Item(331, 2).Value = "seongmo0731@naver.com"
$ws.Cells.Item(331, 3).Value = "경영학과"
$ws.Cells.Item(331, 4).Value = 20192988
$ws.Cells.Item(331, 5).Value = "조성모"
$ws.Cells.Item(331, 6).Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Cells.Item(331, 7).Value = 0.1
$ws.Cells.Item(331, 8).Value = "6:4"
$ws.Cells.Item(331, 9).Value = "10분의 1"
$ws.Cells.Item(331, 10).Value = "20만호, 69만명"
$ws.Cells.Item(331, 11).Value = "충청"
$ws.Cells.Item(331, 12).Value = "Black"
$ws.Cells.Item(331, 14).Value = "모름/무응답"

# Row 332
$ws.Cells.Item(332, 1).Value = 45192.689555590274
$ws.Cells.Item(332, 2).Value = "1004soeun@naver.com"
$ws.Cells.Item(332, 3).Value = "미디어스쿨"
$ws.Cells.Item(332, 4).Value = 20232514
$ws.Cells.Item(332, 5).Value = "김소은"
$ws.Cells.Item(332, 6).Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Cells.Item(332, 7).Value = 0.1
$ws.Cells.Item(332, 8).Value = "6:4"
$ws.Cells.Item(332, 9).Value = "10분의 1"
$ws.Cells.Item(332, 10).Value = "20만호, 69만명"
$ws.Cells.Item(332, 11).Value = "충청"
$ws.Cells.Item(332, 12).Value = "Black"
$ws.Cells.Item(332, 14).Value = "노동자가 과도한 연장근로를 받을 수 있어 반대한다."

# Row 333
$ws.Cells.Item(333, 1).Value = 45192.691171493054
$ws.Cells.Item(333, 2).Value = "kgy5988@naver.com"
$ws.Cells.Item(333, 3).Value = "소프트웨어학부"
$ws.Cells.Item(333, 4).Value = 20203214
$ws.Cells.Item(333, 5).Value = "김진범"
$ws.Cells.Item(333, 6).Value = "‘조(租)’는 공전(公田)의 경작자가 국고에 상납하는 지대 또는 사전(私田)의 경작자가 전주에게 바치는 지대를 뜻한다."
$ws.Cells.Item(333, 7).Value = 0.1
$ws.Cells.Item(333, 8).Value = "6:4"
$ws.Cells.Item(333, 9).Value = "20분의 1"
$ws.Cells.Item(333, 10).Value = "20만호, 69만명"
$ws.Cells.Item(333, 11).Value = "충청"
$ws.Cells.Item(333, 12).Value = "Black"
$ws.Cells.Item(333, 14).Value = "찬성한다."

# Row 334
$ws.Cells.Item(334, 1).Value = 45192.69671621528
$ws.Cells.Item(334, 2).Value = "sysy050300@naver.com"
$ws.Cells.Item(334, 3).Value = "심리학과"
$ws.Cells.Item(334, 4).Value = 20232101
$ws.Cells.Item(334, 5).Value = "고서연"
$ws.Cells.Item(334, 6).Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Cells.Item(334, 7).Value = 0.1
$ws.Cells.Item(334, 8).Value = "6:4"
$ws.Cells.Item(334, 9).Value = "20분의 1"
$ws.Cells.Item(334, 10).Value = "20만호, 69만명"
$ws.Cells.Item(334, 11).Value = "충청"
$ws.Cells.Item(334, 12).Value = "Red"
$ws.Cells.Item(334, 13).Value = "근로시간과 휴무를 유연하게 조정할 수 있어 찬성한다."

# Row 335
$ws.Cells.Item(335, 1).Value = 45192.698860104167
$ws.Cells.Item(335, 2).Value = "ggr1042@naver.com"
$ws.Cells.Item(335, 3).Value = "중국학과"
$ws.Cells.Item(335, 4).Value = 20221542
$ws.Cells.Item(335, 5).Value = "김경록"
$ws.Cells.Item(335, 6).Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Cells.Item(335, 7).Value = 0.5
$ws.Cells.Item(335, 8).Value = "3:7"
$ws.Cells.Item(335, 9).Value = "10분의 1"
$ws.Cells.Item(335, 10).Value = "20만호, 69만명"
$ws.Cells.Item(335, 11).Value = "충청"
$ws.Cells.Item(335, 12).Value = "Red"
$ws.Cells.Item(335, 13).Value = "반대한다."

# Row 336
$ws.Cells.Item(336, 1).Value = 45192.702279074074
$ws.Cells.Item(336, 2).Value = "ntkrud0221@naver.com"
$ws.Cells.Item(336, 3).Value = "체육학과"
$ws.Cells.Item(336, 4).Value = 20234118
$ws.Cells.Item(336, 5).Value = "노태경"
$ws.Cells.Item(336, 6).Value = "‘세(稅)’는 사전의 소유자가 국가에 상납하는 지대를 뜻한다."
$ws.Cells.Item(336, 7).Value = 0.1
$ws.Cells.Item(336, 8).Value = "6:4"
$ws.Cells.Item(336, 9).Value = "15분의 1"
$ws.Cells.Item(336, 10).Value = "20만호, 69만명"
$ws.Cells.Item(336, 11).Value = "경기"
$ws.Cells.Item(336, 12).Value = "Red"
$ws.Cells.Item(336, 13).Value = "모름/무응답"

# Row 337
$ws.Cells.Item(337, 1).Value = 45192.703343784728
$ws.Cells.Item(337, 2).Value = "at79711@naver.com"
$ws.Cells.Item(337, 3).Value = "데이터사이언스"
$ws.Cells.Item(337, 4).Value = 20233220
$ws.Cells.Item(337, 5).Value = "박재영"
$ws.Cells.Item(337, 6).Value = "‘세(稅)’는 사전의 소유자가 국가에 상납하는 지대를 뜻한다."
$ws.Cells.Item(337, 7).Value = 0.5
$ws.Cells.Item(337, 8).Value = "4:6"
$ws.Cells.Item(337, 9).Value = "10분의 1"
$ws.Cells.Item(337, 10).Value = "44만호, 153만명"
$ws.Cells.Item(337, 11).Value = "평안"
$ws.Cells.Item(337, 12).Value = "Black"
$ws.Cells.Item(337, 14).Value = "노동자가 과도한 연장근로를 받을 수 있어 반대한다."

# Row 338
$ws.Cells.Item(338, 1).Value = 45192.712597083329
$ws.Cells.Item(338, 2).Value = "p51008085@gmail.com"
$ws.Cells.Item(338, 3).Value = "인공지능융합학부"
$ws.Cells.Item(338, 4).Value = 20236726
$ws.Cells.Item(338, 5).Value = "박준수"
$ws.Cells.Item(338, 6).Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Cells.Item(338, 7).Value = 0.1
$ws.Cells.Item(338, 8).Value = "6:4"
$ws.Cells.Item(338, 9).Value = "20분의 1"
$ws.Cells.Item(338, 10).Value = "20만호, 69만명"
$ws.Cells.Item(338, 11).Value = "충청"
$ws.Cells.Item(338, 12).Value = "Black"
$ws.Cells.Item(338, 14).Value = "노동자가 과도한 연장근로를 받을 수 있어 반대한다."

# Row 339
$ws.Cells.Item(339, 1).Value = 45192.718088148147
$ws.Cells.Item(339, 2).Value = "jaejae7070@naver.com"
$ws.Cells.Item(339, 3).Value = "일본학과"
$ws.Cells.Item(339, 4).Value = 20221631
$ws.Cells.Item(339, 5).Value = "이재빈"
$ws.Cells.Item(339, 6).Value = "과전법 체제에서 전국 토지를 세 등급으로 나누고 실제 수확량을 확인하여 징수하였다."
$ws.Cells.Item(339, 7).Value = 0.1
$ws.Cells.Item(339, 8).Value = "7:3"
$ws.Cells.Item(339, 9).Value = "10분의 1"
$ws.Cells.Item(339, 10).Value = "15만호,  32만명"
$ws.Cells.Item(339, 11).Value = "경기"
$ws.Cells.Item(339, 12).Value = "Red"
$ws.Cells.Item(339, 13).Value = "반대한다."

# Row 340
$ws.Cells.Item(340, 1).Value = 45192.718404791667
$ws.Cells.Item(340, 2).Value = "leedongbin01@naver.com"
$ws.Cells.Item(340, 3).Value = "영어영문학과"
$ws.Cells.Item(340, 4).Value = 20231224
$ws.Cells.Item(340, 5).Value = "이동빈"
$ws.Cells.Item(340, 6).Value = "‘세(稅)’는 사전의 소유자가 국가에 상납하는 지대를 뜻한다."
$ws.Cells.Item(340, 7).Value = 0.7
$ws.Cells.Item(340, 8).Value = "3:7"
$ws.Cells.Item(340, 9).Value = "10분의 1"
$ws.Cells.Item(340, 10).Value = "20만호, 69만명"
$ws.Cells.Item(340, 11).Value = "평안"
$ws.Cells.Item(340, 12).Value = "Black"
$ws.Cells.Item(340, 14).Value = "찬성한다."

# Row 341
$ws.Cells.Item(341, 1).Value = 45192.726848414357
$ws.Cells.Item(341, 2).Value = "20217134@hallym.ac.kr"
$ws.Cells.Item(341, 3).Value = "체육학과"
$ws.Cells.Item(341, 4).Value = 20217134
$ws.Cells.Item(341, 5).Value = "장효경"
$ws.Cells.Item(341, 6).Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Cells.Item(341, 7).Value = 0.1
$ws.Cells.Item(341, 8).Value = "6:4"
$ws.Cells.Item(341, 9).Value = "10분의 1"
$ws.Cells.Item(341, 10).Value = "20만호, 69만명"
$ws.Cells.Item(341, 11).Value = "평안"
$ws.Cells.Item(341, 12).Value = "Red"
$ws.Cells.Item(341, 13).Value = "근로시간과 휴무를 유연하게 조정할 수 있어 찬성한다."

# Row 342
$ws.Cells.Item(342, 1).Value = 45192.728155578705
$ws.Cells.Item(342, 2).Value = "jinwoo3817@naver.com"
$ws.Cells.Item(342, 3).Value = "디지털미디어콘텐츠"
$ws.Cells.Item(342, 4).Value = 20222552
$ws.Cells.Item(342, 5).Value = "원진우"
$ws.Cells.Item(342, 6).Value = "‘조(租)’는 공전(公田)의 경작자가 국고에 상납하는 지대 또는 사전(私田)의 경작자가 전주에게 바치는 지대를 뜻한다."
$ws.Cells.Item(342, 7).Value = 0.1
$ws.Cells.Item(342, 8).Value = "3:7"
$ws.Cells.Item(342, 9).Value = "10분의 1"
$ws.Cells.Item(342, 10).Value = "20만호, 69만명"
$ws.Cells.Item(342, 11).Value = "충청"
$ws.Cells.Item(342, 12).Value = "Black"
$ws.Cells.Item(342, 14).Value = "모름/무응답"

# Row 343
$ws.Cells.Item(343, 1).Value = 45192.729664895829
$ws.Cells.Item(343, 2).Value = "han7434370@naver.com"
$ws.Cells.Item(343, 3).Value = "체육학과"
$ws.Cells.Item(343, 4).Value = 20224152
$ws.Cells.Item(343, 5).Value = "한진우"
$ws.Cells.Item(343, 6).Value = "‘조(租)’는 공전(公田)의 경작자가 국고에 상납하는 지대 또는 사전(私田)의 경작자가 전주에게 바치는 지대를 뜻한다."
$ws.Cells.Item(343, 7).Value = 0.1
$ws.Cells.Item(343, 8).Value = "4:6"
$ws.Cells.Item(343, 9).Value = "10분의 1"
$ws.Cells.Item(343, 10).Value = "20만호, 69만명"
$ws.Cells.Item(343, 11).Value = "전라"
$ws.Cells.Item(343, 12).Value = "Red"
$ws.Cells.Item(343, 13).Value = "반대한다."

# Row 344
$ws.Cells.Item(344, 1).Value = 45192.735113668983
$ws.Cells.Item(344, 2).Value = "qasw0529@naver.com"
$ws.Cells.Item(344, 3).Value = "미디어스쿨"
$ws.Cells.Item(344, 4).Value = 20232508
$ws.Cells.Item(344, 5).Value = "김민서"
$ws.Cells.Item(344, 6).Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Cells.Item(344, 7).Value = 0.1
$ws.Cells.Item(344, 8).Value = "3:7"
$ws.Cells.Item(344, 9).Value = "10분의 1"
$ws.Cells.Item(344, 10).Value = "20만호, 69만명"
$ws.Cells.Item(344, 11).Value = "전라"
$ws.Cells.Item(344, 12).Value = "Black"
$ws.Cells.Item(344, 14).Value = "찬성한다."

# Row 345
$ws.Cells.Item(345, 1).Value = 45192.740849039357
$ws.Cells.Item(345, 2).Value = "jamesjm0612@gmail.com"
$ws.Cells.Item(345, 3).Value = "영어영문학과"
$ws.Cells.Item(345, 4).Value = 20231231
$ws.Cells.Item(345, 5).Value = "정재민"
$ws.Cells.Item(345, 6).Value = "‘조(租)’는 공전(公田)의 경작자가 국고에 상납하는 지대 또는 사전(私田)의 경작자가 전주에게 바치는 지대를 뜻한다."
$ws.Cells.Item(345, 7).Value = 0.3
$ws.Cells.Item(345, 8).Value = "6:4"
$ws.Cells.Item(345, 9).Value = "15분의 1"
$ws.Cells.Item(345, 10).Value = "20만호, 69만명"
$ws.Cells.Item(345, 11).Value = "평안"
$ws.Cells.Item(345, 12).Value = "Black"
$ws.Cells.Item(345, 14).Value = "모름/무응답"

# Row 346
$ws.Cells.Item(346, 1).Value = 45192.749265208331
$ws.Cells.Item(346, 2).Value = "meldek98@gmail.com"
$ws.Cells.Item(346, 3).Value = "소프트웨어학부"
$ws.Cells.Item(346, 4).Value = 20235102
$ws.Cells.Item(346, 5).Value = "강비성"
$ws.Cells.Item(346, 6).Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Cells.Item(346, 7).Value = 0.1
$ws.Cells.Item(346, 8).Value = "6:4"
$ws.Cells.Item(346, 9).Value = "20분의 1"
$ws.Cells.Item(346, 10).Value = "20만호, 69만명"
$ws.Cells.Item(346, 11).Value = "충청"
$ws.Cells.Item(346, 12).Value = "Black"
$ws.Cells.Item(346, 14).Value = "찬성한다."

# Row 347
$ws.Cells.Item(347, 1).Value = 45192.760950810189
$ws.Cells.Item(347, 2).Value = "alscoco100@gmail.com"
$ws.Cells.Item(347, 3).Value = "식품영양학과"
$ws.Cells.Item(347, 4).Value = 20223806
$ws.Cells.Item(347, 5).Value = "김민채"
$ws.Cells.Item(347, 6).Value = "‘조(租)’는 공전(公田)의 경작자가 국고에 상납하는 지대 또는 사전(私田)의 경작자가 전주에게 바치는 지대를 뜻한다."
$ws.Cells.Item(347, 7).Value = 0.7
$ws.Cells.Item(347, 8).Value = "5:5"
$ws.Cells.Item(347, 9).Value = "20분의 1"
$ws.Cells.Item(347, 10).Value = "15만호,  32만명"
$ws.Cells.Item(347, 11).Value = "충청"
$ws.Cells.Item(347, 12).Value = "Black"
$ws.Cells.Item(347, 14).Value = "모름/무응답"

# Row 348
$ws.Cells.Item(348, 1).Value = 45192.767777824076
$ws.Cells.Item(348, 2).Value = "yeshin05@naver.com"
$ws.Cells.Item(348, 3).Value = "미래융합스쿨"
$ws.Cells.Item(348, 4).Value = 20236639
$ws.Cells.Item(348, 5).Value = "최예원"
$ws.Cells.Item(348, 6).Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Cells.Item(348, 7).Value = 0.1
$ws.Cells.Item(348, 8).Value = "6:4"
$ws.Cells.Item(348, 9).Value = "20분의 1"
$ws.Cells.Item(348, 10).Value = "20만호, 69만명"
$ws.Cells.Item(348, 11).Value = "충청"
$ws.Cells.Item(348, 12).Value = "Red"
$ws.Cells.Item(348, 13).Value = "근로시간과 휴무를 유연하게 조정할 수 있어 찬성한다."

# Row 349
$ws.Cells.Item(349, 1).Value = 45192.773631249998
$ws.Cells.Item(349, 2).Value = "eugene3551@gmail.com"
$ws.Cells.Item(349, 3).Value = "소프트웨어학부"
$ws.Cells.Item(349, 4).Value = 20235214
$ws.Cells.Item(349, 5).Value = "유수영"
$ws.Cells.Item(349, 6).Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Cells.Item(349, 7).Value = 0.9
$ws.Cells.Item(349, 8).Value = "6:4"
$ws.Cells.Item(349, 9).Value = "30분의 1"
$ws.Cells.Item(349, 10).Value = "20만호, 69만명"
$ws.Cells.Item(349, 11).Value = "충청"
$ws.Cells.Item(349, 12).Value = "Red"
$ws.Cells.Item(349, 13).Value = "반대한다."

# Row 350
$ws.Cells.Item(350, 1).Value = 45192.78515171296
$ws.Cells.Item(350, 2).Value = "cozyandrelaxing2@gmail.com"
$ws.Cells.Item(350, 3).Value = "환경생명공학과"
$ws.Cells.Item(350, 4).Value = 20223725
$ws.Cells.Item(350, 5).Value = "이성민"
$ws.Cells.Item(350, 6).Value = "실제로 현장에 나가서 수확량을 파악하고 등급을 매기는 답험(踏驗)을 하였다."
$ws.Cells.Item(350, 7).Value = 0.7
$ws.Cells.Item(350, 8).Value = "4:6"
$ws.Cells.Item(350, 9).Value = "10분의 1"
$ws.Cells.Item(350, 10).Value = "20만호, 69만명"
$ws.Cells.Item(350, 11).Value = "전라"
$ws.Cells.Item(350, 12).Value = "Black"
$ws.Cells.Item(350, 14).Value = "모름/무응답"

# Row 351
$ws.Cells.Item(351, 1).Value = 45192.818970231485
$ws.Cells.Item(351, 2).Value = "hyunbin7379@gmail.com"
$ws.Cells.Item(351, 3).Value = "경영학부"
$ws.Cells.Item(351, 4).Value = 20233036
$ws.Cells.Item(351, 5).Value = "정현빈"
$ws.Cells.Item(351, 6).Value = "과전법 체제에서 전국 토지를 세 등급으로 나누고 실제 수확량을 확인하여 징수하였다."
$ws.Cells.Item(351, 7).Value = 0.5
$ws.Cells.Item(351, 8).Value = "5:5"
$ws.Cells.Item(351, 9).Value = "15분의 1"
$ws.Cells.Item(351, 10).Value = "130만호, 5백만명"
$ws.Cells.Item(351, 11).Value = "경기"
$ws.Cells.Item(351, 12).Value = "Red"
$ws.Cells.Item(351, 13).Value = "모름/무응답"

# Row 352
$ws.Cells.Item(352, 1).Value = 45192.82088186343
$ws.Cells.Item(352, 2).Value = "jyn10131@naver.com"
$ws.Cells.Item(352, 3).Value = "식품영양학과"
$ws.Cells.Item(352, 4).Value = 20233847
$ws.Cells.Item(352, 5).Value = "정예나"
$ws.Cells.Item(352, 6).Value = "‘조(租)’는 공전(公田)의 경작자가 국고에 상납하는 지대 또는 사전(私田)의 경작자가 전주에게 바치는 지대를 뜻한다."
$ws.Cells.Item(352, 7).Value = 0.3
$ws.Cells.Item(352, 8).Value = "6:4"
$ws.Cells.Item(352, 9).Value = "15분의 1"
$ws.Cells.Item(352, 10).Value = "44만호, 153만명"
$ws.Cells.Item(352, 11).Value = "경상"
$ws.Cells.Item(352, 12).Value = "Red"
$ws.Cells.Item(352, 13).Value = "근로시간과 휴무를 유연하게 조정할 수 있어 찬성한다."

# 4. Update view state: freeze header row, scroll so row 257 is top-left of the
#    scrollable pane, and select D274 (matches the final selection in the sheet).
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$excel.ActiveWindow.ScrollRow = 257
$ws.Range("D274").Select()

$excel.ScreenUpdating = $true